$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 3429
$ws.Range("J3").Value = 3593
$ws.Range("J4").Value = 796
$ws.Range("J5").Value = 281
$ws.Range("J6").Value = 4197
$ws.Range("J7").Value = 12296

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J7").Value = 376
$ws.Range("J8").Value = 806
$ws.Range("J11").Value = 189
$ws.Range("J17").Value = 21
$ws.Range("J19").Value = 375
$ws.Range("J20").Value = 257
$ws.Range("J21").Value = 22
$ws.Range("J27").Value = 75
$ws.Range("J29").Value = 698
$ws.Range("J31").Value = 96
$ws.Range("J33").Value = 559
$ws.Range("J34").Value = 62
$ws.Range("J36").Value = 181
$ws.Range("J37").Value = 391
$ws.Range("J42").Value = 473
$ws.Range("J47").Value = 89
$ws.Range("J48").Value = 124
$ws.Range("J49").Value = 81
$ws.Range("J51").Value = 164
$ws.Range("J52").Value = 339
$ws.Range("J55").Value = 151
$ws.Range("J61").Value = 19
$ws.Range("J63").Value = 54
$ws.Range("J65").Value = 323
$ws.Range("J67").Value = 452
$ws.Range("J68").Value = 20
$ws.Range("J73").Value = 106
$ws.Range("J75").Value = 38
$ws.Range("J78").Value = 162
$ws.Range("J79").Value = 365
$ws.Range("J83").Value = 288
$ws.Range("J84").Value = 109
$ws.Range("J85").Value = 559
$ws.Range("J86").Value = 68
$ws.Range("J87").Value = 40
$ws.Range("J88").Value = 132
$ws.Range("J89").Value = 145
$ws.Range("J95").Value = 191
$ws.Range("J96").Value = 134
$ws.Range("J97").Value = 77
$ws.Range("J99").Value = 176
$ws.Range("J101").Value = 12296

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 139
$ws.Range("J3").Value = 210
$ws.Range("J4").Value = 43
$ws.Range("J6").Value = 156
$ws.Range("J7").Value = 559

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J6").Value = 145
$ws.Range("J7").Value = 339

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J2").Value = 69
$ws.Range("J7").Value = 189

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 237
$ws.Range("J3").Value = 253
$ws.Range("J7").Value = 806

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J6").Value = 125
$ws.Range("J7").Value = 376

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J3").Value = 38
$ws.Range("J4").Value = 18
$ws.Range("J7").Value = 145

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J2").Value = 40
$ws.Range("J7").Value = 134

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 118
$ws.Range("J3").Value = 137
$ws.Range("J4").Value = 9
$ws.Range("J6").Value = 111
$ws.Range("J7").Value = 391

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J3").Value = 59
$ws.Range("J6").Value = 48
$ws.Range("J7").Value = 176

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J3").Value = 191
$ws.Range("J7").Value = 452

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J3").Value = 24
$ws.Range("J6").Value = 25
$ws.Range("J7").Value = 96

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J2").Value = 36
$ws.Range("J7").Value = 109

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J3").Value = 97
$ws.Range("J6").Value = 115
$ws.Range("J7").Value = 323

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 84
$ws.Range("J3").Value = 107
$ws.Range("J6").Value = 82
$ws.Range("J7").Value = 288

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J3").Value = 58
$ws.Range("J7").Value = 191

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J5").Value = 24
$ws.Range("J6").Value = 180
$ws.Range("J7").Value = 559

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J4").Value = 6
$ws.Range("J6").Value = 46
$ws.Range("J7").Value = 81

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 215
$ws.Range("J3").Value = 241
$ws.Range("J7").Value = 698

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 90
$ws.Range("J3").Value = 106
$ws.Range("J6").Value = 142
$ws.Range("J7").Value = 375

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J6").Value = 58
$ws.Range("J7").Value = 124

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 99
$ws.Range("J5").Value = 13
$ws.Range("J7").Value = 473

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J2").Value = 41
$ws.Range("J3").Value = 55
$ws.Range("J7").Value = 162

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J3").Value = 33
$ws.Range("J7").Value = 151

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("J6").Value = 14
$ws.Range("J7").Value = 22

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J3").Value = 137
$ws.Range("J7").Value = 365

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J3").Value = 82
$ws.Range("J6").Value = 64
$ws.Range("J7").Value = 257

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("J3").Value = 9
$ws.Range("J7").Value = 21

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J2").Value = 66
$ws.Range("J4").Value = 4
$ws.Range("J7").Value = 181

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 62

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 89

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 106

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J2").Value = 15
$ws.Range("J7").Value = 77

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J3").Value = 43
$ws.Range("J6").Value = 51
$ws.Range("J7").Value = 132

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J4").Value = 9
$ws.Range("J7").Value = 75

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J4").Value = 32
$ws.Range("J7").Value = 68

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 38

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J2").Value = 42
$ws.Range("J3").Value = 44
$ws.Range("J6").Value = 53
$ws.Range("J7").Value = 164

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("J2").Value = 8
$ws.Range("J7").Value = 20

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("J6").Value = 22
$ws.Range("J7").Value = 40

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range("J6").Value = 5
$ws.Range("J7").Value = 19
